$wb = $excel.ActiveWorkbook

# --- Sheets ---------------------------------------------------------------
$wsCenit = $wb.Worksheets.Item(1)   # "Transporte Cenit"
$wsTGI   = $wb.Worksheets.Item(2)   # "TGI Transporte de Gas"
$wsProd  = $wb.Worksheets.Item(3)   # "Producción"

# --- sheet2 (TGI Transporte de Gas): fix the shared-formula range ---------
# E13 was the anchor of the shared formula ref="E9:E19"; it should only cover
# E13:E19 (E9:E12 keep their own individual formulas).
$wsTGI.Range("E13").Formula = "=+D13/C13"
$wsTGI.Range("E14:E19").FormulaR1C1 = "=+RC[-1]/RC[-2]"

# --- sheet3 (Producción): update the data table ----------------------------
$wsProd.Range("H7").Value2 = 0.52
$wsProd.Range("I7").Value = "Kwh/BE"

$wsProd.Range("G10").Value = "Total produccion"
$wsProd.Range("I10").Value2 = 750
$wsProd.Range("J10").Value = "kpd"

$wsProd.Range("I11").Formula = "=+I10*1000*365"
$wsProd.Range("J11").Value = "barriles/año"

$wsProd.Range("G13").Value = "Electricidad T"
$wsProd.Range("I13").Formula = "=+I11*H7"
$wsProd.Range("J13").Value = "Kwh"

$wsProd.Range("I14").NumberFormat = "0.0"
$wsProd.Range("I14").Formula = "=+I13/1000000"
$wsProd.Range("J14").Value = "Gwh"

# --- selections -------------------------------------------------------------
$wsCenit.Range("F22").Select()
$wsTGI.Range("E40").Select()
$wsProd.Range("H8").Select()

# --- active sheet / tab -----------------------------------------------------
# Activating "Producción" last makes it the active tab and moves
# tabSelected="1" off "TGI Transporte de Gas" onto it, matching the diff
# (activeTab 1 -> 2).
$wsProd.Activate()
